# Add the missing "temp_june" (column G) readings for every data row
# (rows 2-73 on sheet "data"). The header G1 ("temp_june") already exists;
# only the per-plot values were missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

$juneTemps = @(
    26.1, 26.1, 26.3, 26.5, 26.6, 26.9, 26.7, 26.6, 26.6, 26.1,
    25.8, 26.6, 26.5, 26.0, 26.7, 26.7, 26.9, 26.6, 27.3, 27.2,
    27.2, 26.6, 26.8, 27.2, 27.5, 27.7, 27.7, 27.8, 27.2, 27.6,
    27.0, 27.0, 27.3, 26.5, 26.7, 26.7, 28.2, 28.3, 28.3, 28.7,
    28.4, 28.6, 28.5, 28.9, 28.8, 27.7, 28.8, 28.8, 27.9, 27.3,
    28.3, 28.3, 28.0, 28.5, 28.7, 27.2, 28.3, 27.6, 27.9, 28.4,
    28.3, 28.3, 27.8, 27.6, 27.1, 27.4, 28.3, 28.1, 27.6, 27.5,
    27.1, 26.7
)

$firstRow = 2
for ($i = 0; $i -lt $juneTemps.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $juneTemps[$i]
}

# Restore the view state recorded in the saved workbook: scrolled so row 52
# is at the top and the last-entered cell (G74, just below the data) selected.
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G74").Select()
